{"js": "// Fix typo \"Anergy\" -> \"Energy\" in the company name \"Legend Anergy Advisors\"\n// (commit message: \"update resume - typo\").\nconst body = context.document.body;\n\n// Locate the exact run text that contains the typo.\nconst results = body.search(\"Legend Anergy Advisors\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'Legend Anergy Advisors' in the document.\");\n}\n\n// Replace in place so the surrounding run formatting (Arial, bold, size 20,\n// black) is preserved; this also keeps the trailing space that followed\n// \"Advisors\" in the original run untouched.\nresults.items[0].insertText(\"Legend Energy Advisors\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Fix typo \"Anergy\" -> \"Energy\" in the company name \"Legend Anergy Advisors\"\n# (commit message: \"update resume - typo\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Legend Anergy Advisors\"\n$find.Replacement.Text = \"Legend Energy Advisors\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace(wdReplaceAll=2)\n$find.Execute(\"Legend Anergy Advisors\", $true, $false, $false, $false, $false, $true, 1, $false, \"Legend Energy Advisors\", 2)\n"}
